$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / shared-string updates ---
# Mayor name (M6, merged M6:N6)
$ws.Range("M6").Value = "Thomas G. Donlon"

# Volume / Number line (A8, merged A8:B8... actually C8:L8 per merges, but shared string is on A8)
$ws.Range("A8").Value = "Volume 31   Number  39"

# Report covering the week line (C9)
$ws.Range("C9").Value = "Report Covering the Week  9/23/2024  Through  9/29/2024"

# --- Fix up cell types/styles that change between numeric and text (copy format+value from a
#     donor cell that already has the exact target style + shared text, or target style for numbers) ---
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("C14").Copy($ws.Range("C30"))
$ws.Range("C14").Copy($ws.Range("D33"))
$ws.Range("E14").Copy($ws.Range("E33"))

$ws.Range("C15").Copy($ws.Range("D31"))
$ws.Range("E15").Copy($ws.Range("E31"))
$ws.Range("C15").Copy($ws.Range("F31"))

# --- Numeric value updates across the weekly crime-stat table (rows 14-33) ---
$ws.Range("M14").Value = 45.454545454545
$ws.Range("N14").Value = -75
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 3
$ws.Range("F15").Value = 18
$ws.Range("G15").Value = 11
$ws.Range("H15").Value = 63.636363636363
$ws.Range("I15").Value = 119
$ws.Range("J15").Value = 88
$ws.Range("K15").Value = 35.227272727272
$ws.Range("L15").Value = -10.526315789473
$ws.Range("M15").Value = 60.810810810810
$ws.Range("N15").Value = -30.813953488372
$ws.Range("C16").Value = 23
$ws.Range("D16").Value = 28
$ws.Range("E16").Value = -17.857142857142
$ws.Range("F16").Value = 126
$ws.Range("G16").Value = 136
$ws.Range("H16").Value = -7.352941176470
$ws.Range("I16").Value = 1242
$ws.Range("J16").Value = 1353
$ws.Range("K16").Value = -8.203991130820
$ws.Range("L16").Value = -21.541377132027
$ws.Range("M16").Value = 31.567796610169
$ws.Range("N16").Value = -84.304309364337
$ws.Range("C17").Value = 41
$ws.Range("D17").Value = 44
$ws.Range("E17").Value = -6.818181818181
$ws.Range("F17").Value = 179
$ws.Range("G17").Value = 167
$ws.Range("H17").Value = 7.185628742514
$ws.Range("I17").Value = 1755
$ws.Range("J17").Value = 1655
$ws.Range("K17").Value = 6.042296072507
$ws.Range("L17").Value = 10.031347962382
$ws.Range("M17").Value = 88.102893890675
$ws.Range("N17").Value = -32.265534542647
$ws.Range("C18").Value = 28
$ws.Range("D18").Value = 35
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 124
$ws.Range("G18").Value = 154
$ws.Range("H18").Value = -19.480519480519
$ws.Range("I18").Value = 1365
$ws.Range("J18").Value = 1578
$ws.Range("K18").Value = -13.498098859315
$ws.Range("L18").Value = -39.601769911504
$ws.Range("M18").Value = 1.865671641791
$ws.Range("N18").Value = -84.753713838936
$ws.Range("C19").Value = 236
$ws.Range("D19").Value = 224
$ws.Range("E19").Value = 5.357142857142
$ws.Range("F19").Value = 866
$ws.Range("G19").Value = 939
$ws.Range("H19").Value = -7.774227902023
$ws.Range("I19").Value = 7773
$ws.Range("J19").Value = 8633
$ws.Range("K19").Value = -9.961774585891
$ws.Range("L19").Value = -9.888708555529
$ws.Range("M19").Value = 0.504266873545
$ws.Range("N19").Value = -68.317437026167
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = -18.181818181818
$ws.Range("F20").Value = 43
$ws.Range("G20").Value = 71
$ws.Range("H20").Value = -39.436619718309
$ws.Range("I20").Value = 343
$ws.Range("J20").Value = 481
$ws.Range("K20").Value = -28.690228690228
$ws.Range("L20").Value = -33.783783783783
$ws.Range("M20").Value = 13.576158940397
$ws.Range("N20").Value = -92.517452006980
$ws.Range("C21").Value = 339
$ws.Range("D21").Value = 345
$ws.Range("E21").Value = -1.739130434782
$ws.Range("F21").Value = 1357
$ws.Range("G21").Value = 1478
$ws.Range("H21").Value = -8.186738836265
$ws.Range("I21").Value = 12613
$ws.Range("J21").Value = 13802
$ws.Range("K21").Value = -8.614693522677
$ws.Range("L21").Value = -14.435927006308
$ws.Range("M21").Value = 11.245369553713
$ws.Range("N21").Value = -74.159513224478
$ws.Range("C22").Value = 13
$ws.Range("D22").Value = 9
$ws.Range("E22").Value = 44.444444444444
$ws.Range("F22").Value = 41
$ws.Range("G22").Value = 42
$ws.Range("H22").Value = -2.380952380952
$ws.Range("I22").Value = 453
$ws.Range("J22").Value = 486
$ws.Range("K22").Value = -6.790123456790
$ws.Range("L22").Value = -8.853118712273
$ws.Range("M22").Value = 14.105793450881
$ws.Range("C23").Value = 4
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 26
$ws.Range("G23").Value = 36
$ws.Range("H23").Value = -27.777777777777
$ws.Range("I23").Value = 293
$ws.Range("J23").Value = 305
$ws.Range("K23").Value = -3.934426229508
$ws.Range("L23").Value = -15.317919075144
$ws.Range("M23").Value = 11.406844106463
$ws.Range("C24").Value = 487
$ws.Range("D24").Value = 413
$ws.Range("E24").Value = 17.917675544794
$ws.Range("F24").Value = 1860
$ws.Range("G24").Value = 1583
$ws.Range("H24").Value = 17.498420720151
$ws.Range("I24").Value = 16664
$ws.Range("J24").Value = 15473
$ws.Range("K24").Value = 7.697279131390
$ws.Range("L24").Value = 1.746244962754
$ws.Range("M24").Value = 32.065303534633
$ws.Range("C25").Value = 408
$ws.Range("D25").Value = 342
$ws.Range("E25").Value = 19.298245614035
$ws.Range("F25").Value = 1511
$ws.Range("G25").Value = 1270
$ws.Range("H25").Value = 18.976377952755
$ws.Range("I25").Value = 14059
$ws.Range("J25").Value = 12816
$ws.Range("K25").Value = 9.698813982521
$ws.Range("L25").Value = 0.918814155480
$ws.Range("C26").Value = 110
$ws.Range("D26").Value = 90
$ws.Range("E26").Value = 22.222222222222
$ws.Range("F26").Value = 455
$ws.Range("G26").Value = 375
$ws.Range("H26").Value = 21.333333333333
$ws.Range("I26").Value = 3801
$ws.Range("J26").Value = 3737
$ws.Range("K26").Value = 1.712603692801
$ws.Range("L26").Value = 8.879977083930
$ws.Range("M26").Value = 40.361890694239
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 16
$ws.Range("H27").Value = 56.25
$ws.Range("I27").Value = 181
$ws.Range("J27").Value = 156
$ws.Range("K27").Value = 16.025641025641
$ws.Range("L27").Value = -13.397129186602
$ws.Range("C28").Value = 21
$ws.Range("D28").Value = 23
$ws.Range("E28").Value = -8.695652173913
$ws.Range("F28").Value = 89
$ws.Range("G28").Value = 66
$ws.Range("H28").Value = 34.848484848484
$ws.Range("I28").Value = 729
$ws.Range("J28").Value = 677
$ws.Range("K28").Value = 7.680945347119
$ws.Range("L28").Value = -0.409836065573
$ws.Range("F29").Value = 5
$ws.Range("L29").Value = -26.829268292682
$ws.Range("M29").Value = 7.142857142857
$ws.Range("F30").Value = 5
$ws.Range("L30").Value = -27.777777777777
$ws.Range("M30").Value = 18.181818181818
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 9
$ws.Range("H31").Value = -10
$ws.Range("I31").Value = 113
$ws.Range("J31").Value = 97
$ws.Range("K31").Value = 16.494845360824
$ws.Range("L31").Value = -11.71875
$ws.Range("I33").Value = 17
$ws.Range("K33").Value = 41.666666666666
$ws.Range("L33").Value = 54.545454545454
